$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C:\Users\Pascal.Kipkemoi\git\MKOPA_Regression_Test_Channel\FilesToUpload\Mpesalatest File.csv"
$ws.Range("C2").Value = "polr780"
